# Repull data, push all data, mean calculation
# Updates the "dSF" column (F) values for the quantrill_cal calendar sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new dSF (column F) value, as scraped/recalculated from the
# refreshed source data.
$updates = @{
    2  = -5
    4  = -2
    5  = 3
    6  = -2
    7  = 3
    8  = -2
    9  = 1
    10 = 10
    11 = -3
    12 = -5
    13 = 1
    15 = 1
    16 = -2
    17 = 4
    18 = 3
    19 = -2
    20 = 3
    21 = 1
    22 = 6
    23 = 6
    24 = 3
    25 = -2
    28 = -2
    29 = -4
    30 = -2
    31 = -4
    32 = 4
    33 = 2
    34 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
